# Add a "Save" column (column H) to the s_vals worksheet.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell, same bold/bordered style as the other header cells (B1:G1)
$ws.Range("H1").Value = "Save"
$ws.Range("G1").Copy()
$ws.Range("H1").PasteSpecial(-4122)

# "Save" values for each data row (1 = win saved, 0 = not)
$saveValues = @(1,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,1,0,0)

for ($i = 0; $i -lt $saveValues.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 8).Value = $saveValues[$i]
}
